$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'running leggings for women high waist',
    'running leggings for women pack',
    'running leggings for women winter',
    'running leggings pack',
    'running leggings pants',
    'running leggings winter women',
    'running leggings women',
    'running legs men',
    'running light',
    'running light women',
    'running pamts men',
    'running pant',
    'running pant women',
    'running pants',
    'running pants cold',
    'running pants cold weather',
    'running pants cold weather men',
    'running pants cold weather womens',
    'running pants compression',
    'running pants compression men',
    'running pants for women',
    'running pants for women cold weather',
    'running pants high waist',
    'running pants leggings',
    'running pants long',
    'running pants loose',
    'running pants men cold',
    'running pants men cold weather',
    'running pants men tights',
    'running pants on',
    'running pants pack',
    'running pants thin',
    'running pants tights',
    'running pants with',
    'running pants with iphone pocket',
    'running pants women',
    'running pants women cold weather',
    'running pants women compression',
    'running pants women for cold weather',
    'running pants women leggings',
    'running pants women reflective',
    'running pants women tall',
    'running pants women winter',
    'running pants womens',
    'running pantyhose',
    'running recovery',
    'running recovery pants',
    'running recovery supplement',
    'running reflective',
    'running reflective clothing',
    'running reflective gear women',
    'running reflective jacket women',
    'running reflective wear',
    'running resistance',
    'running short tights women',
    'running shorts',
    'running shorts compression',
    'running shorts compression women',
    'running shorts for women pack',
    'running shorts for women with underwear',
    'running shorts long women',
    'running shorts tight',
    'running shorts tights',
    'running shorts with compression',
    'running shorts woman',
    'running shorts women',
    'running shorts women compression',
    'running shorts women lululemon',
    'running shorts women pack',
    'running skirt with tights',
    'running skirts with leggings for women',
    'running spandex',
    'running stocking',
    'running supplements for endurance',
    'running support pants',
    'running tape knee',
    'running target',
    'running tighs',
    'running tight women',
    'running tights',
    'running tights capri',
    'running tights cold weather',
    'running tights compression',
    'running tights for men',
    'running tights for women',
    'running tights for women cold weather',
    'running tights high waist',
    'running tights men 2xu',
    'running tights men 3 4',
    'running tights men pocket',
    'running tights men tall',
    'running tights pocket',
    'running tights tall women',
    'running tights thermal',
    'running tights with drawstring',
    'running tights with drawstring women',
    'running tights women',
    'running tights women capri',
    'running tights women cold',
    'running tights women cold weather'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
